$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 82-96 currently only have column B (Concept Id) populated.
# Fill in column C (Default values) with "OPB" and column D
# (Generation type) with "String_test" for each of them, matching every
# other "Diagnostik" row already fully populated further up the sheet.
for ($r = 82; $r -le 96; $r++) {
    $ws.Cells.Item($r, 3).Value = "OPB"
    $ws.Cells.Item($r, 4).Value = "String_test"
}

# Move the active selection to G86, matching the saved view state.
$ws.Range("G86").Select()
